# Corrección a Diebold Mariano y revisión de Cap1
# Update the p-values matrix (P_valores) and the DM statistics matrix (Estadisticos_DM)

$wb = $excel.ActiveWorkbook

# --- Sheet: P_valores ---
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.05229962507050767
$wsP.Range("D2").Value = 0.9607996369289524
$wsP.Range("E2").Value = 0.2755592418626311
$wsP.Range("F2").Value = 0.1076664472397502

$wsP.Range("B3").Value = 0.05229962507050767
$wsP.Range("D3").Value = 0.05174229652590956
$wsP.Range("E3").Value = 0.1925068256742337
$wsP.Range("F3").Value = 0.4780249664159233

$wsP.Range("B4").Value = 0.9607996369289524
$wsP.Range("C4").Value = 0.05174229652590956
$wsP.Range("E4").Value = 0.4962014189058628
$wsP.Range("F4").Value = 0.06702413307398336

$wsP.Range("B5").Value = 0.2755592418626311
$wsP.Range("C5").Value = 0.1925068256742337
$wsP.Range("D5").Value = 0.4962014189058628
$wsP.Range("F5").Value = 0.217733477186375

$wsP.Range("B6").Value = 0.1076664472397502
$wsP.Range("C6").Value = 0.4780249664159233
$wsP.Range("D6").Value = 0.06702413307398336
$wsP.Range("E6").Value = 0.217733477186375

# --- Sheet: Estadisticos_DM ---
$wsE = $wb.Worksheets.Item("Estadisticos_DM")

$wsE.Range("C2").Value = -2.051604956423864
$wsE.Range("D2").Value = 0.04971305756704497
$wsE.Range("E2").Value = -1.118162711753777
$wsE.Range("F2").Value = -1.677148356360227

$wsE.Range("B3").Value = 2.051604956423864
$wsE.Range("D3").Value = 2.056920129174218
$wsE.Range("E3").Value = 1.344436399022186
$wsE.Range("F3").Value = -0.7217911936876183

$wsE.Range("B4").Value = -0.04971305756704497
$wsE.Range("C4").Value = -2.056920129174218
$wsE.Range("E4").Value = -0.691959674395921
$wsE.Range("F4").Value = -1.926760552702913

$wsE.Range("B5").Value = 1.118162711753777
$wsE.Range("C5").Value = -1.344436399022186
$wsE.Range("D5").Value = 0.691959674395921
$wsE.Range("F5").Value = -1.268913754878327

$wsE.Range("B6").Value = 1.677148356360227
$wsE.Range("C6").Value = 0.7217911936876183
$wsE.Range("D6").Value = 1.926760552702913
$wsE.Range("E6").Value = 1.268913754878327

$wb.Save()
